# "break out stock.yaml completed"
# - Convert D67:D76 (bsecode) on the "day" sheet from text to numeric values.
# - Append 8 new rows (77-84) of "day" stock data pulled in on 03/07/2024.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("day")

# --- Fix D67:D76 so bsecode is stored as a real number, not text ---
$bsecodes = @{
    67 = 532977
    68 = 500002
    69 = 500034
    70 = 500182
    71 = 505200
    72 = 500114
    73 = 500520
    74 = 532868
    75 = 535755
    76 = 500470
}
foreach ($row in $bsecodes.Keys) {
    $ws.Cells.Item($row, 4).Value = $bsecodes[$row]
}

# --- Append new rows 77-84 ---
# columns: A=sr B=nsecode C=name D=bsecode E=per_chg F=close G=volume H=timeframe I=Date Time
$newRows = @(
    @(1, "HAL",        "Hindustan Aeronautics Ltd",                            "541154", 2.16,               5459.3,  2852572, "day", "03/07/2024 11:34:46"),
    @(2, "BALKRISIND",  "Balkrishna Industries Limited",                       "502355", -0.6899999999999999, 3108.05, 350922,  "day", "03/07/2024 11:34:46"),
    @(3, "TORNTPHARM",  "Torrent Pharmaceuticals Limited",                     "500420", 0.63,               2837.5,  218267,  "day", "03/07/2024 11:34:46"),
    @(4, "CIPLA",       "Cipla Limited",                                       "500087", -0.33,              1483.75, 1061247, "day", "03/07/2024 11:34:46"),
    @(5, "INDUSINDBK",  "Indusind Bank Limited",                               "532187", 1.82,               1455.5,  5907245, "day", "03/07/2024 11:34:46"),
    @(6, "VOLTAS",      "Voltas Limited",                                      "500575", 0.12,               1451.85, 1022098, "day", "03/07/2024 11:34:46"),
    @(7, "CHOLAFIN",    "Cholamandalam Investment And Finance Company Limited","511243", 2.5,                1435.3,  525587,  "day", "03/07/2024 11:34:46"),
    @(8, "HDFCLIFE",    "HDFC Life Insurance Company Ltd",                     "540777", 1.08,               596.15,  2134780, "day", "03/07/2024 11:34:46")
)

$startRow = 77
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]

    # bsecode stays text for the new rows, matching the source feed
    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = $data[3]

    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
    $ws.Cells.Item($r, 9).Value = $data[8]
}

Write-Output "applied bsecode numeric fix + appended rows 77-84"
